# Release V3.2: Menu Interativo, Filtro de Qualificacao e Visual Limpo
#
# - Rename the roster entries in row 2 (Sai_Nome / Entra_Nome columns)
# - Make the header row font color an explicit black (RGB) instead of the
#   theme-based color, for a "cleaner" visual look

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two name cells on the data row
$ws.Range("B2").Value = "Cb João"
$ws.Range("C2").Value = "Cb Artur"

# Visual Limpo: force the header row text color to explicit black
$ws.Range("A1:C1").Font.Color = 0
